$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: cohort 2020 (A7 index), period 5 -> num_customers 32 -> 33; retention_rate recalculated
$ws.Range("C22").Value = 33
$ws.Range("E22").Value = 33 / 2654

# Row 27: cohort 2021, period 4 -> num_customers 57 -> 58; retention_rate recalculated
$ws.Range("C27").Value = 58
$ws.Range("E27").Value = 58 / 2252

# Row 31: cohort 2022, period 3 -> num_customers 58 -> 61; retention_rate recalculated
$ws.Range("C31").Value = 61
$ws.Range("E31").Value = 61 / 2312

# Row 36: cohort 2024, period 1 -> num_customers 147 -> 148; retention_rate recalculated
$ws.Range("C36").Value = 148
$ws.Range("E36").Value = 148 / 1930

# Row 37: cohort 2025, period 0 -> num_customers and cohort_size 982 -> 1003; retention_rate stays 1
$ws.Range("C37").Value = 1003
$ws.Range("D37").Value = 1003
$ws.Range("E37").Value = 1
